$d = $word.ActiveDocument
$t = $d.Tables.Item(2)

$t.Cell(3, 3).Range.Text = "True"
$t.Cell(3, 4).Range.Text = "True"

$t.Cell(4, 3).Range.Text = "False"
$t.Cell(4, 4).Range.Text = "False"

$t.Cell(5, 3).Range.Text = "False"
$t.Cell(5, 4).Range.Text = "False"
